$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster becomes "ECs") ---
$ws.Range("A2").Value = "ECs"
$ws.Range("G2").Value = 0.03814
$ws.Range("H2").Value = 0.11442
$ws.Range("I2").Value = 0.0004360684493923871
$ws.Range("J2").Value = 0.0004360684493923871
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 0.01292722245333333
$ws.Range("R2").Value = 0.11634500208
$ws.Range("S2").Value = 0.0004360684493923871
$ws.Range("T2").Value = 0.0004360684493923871

# --- Row 3 (Sending cluster becomes "FAPs", was "MuSCs") ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2972863333333333
$ws.Range("H3").Value = 0.891859
$ws.Range("I3").Value = 0.003398982443686812
$ws.Range("J3").Value = 0.003398982443686811
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 0.1007626262017778
$ws.Range("R3").Value = 0.9068636358159999
$ws.Range("S3").Value = 0.003398982443686812
$ws.Range("T3").Value = 0.003398982443686811

# --- Row 4 (Sending cluster stays "Resolving-Mac") ---
$ws.Range("G4").Value = 87.12790666666666
$ws.Range("H4").Value = 261.38372
$ws.Range("I4").Value = 0.9961649491069209
$ws.Range("J4").Value = 0.9961649491069208
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3389413333333333
$ws.Range("N4").Value = 1.016824
$ws.Range("Q4").Value = 29.53124885614222
$ws.Range("R4").Value = 265.78123970528
$ws.Range("S4").Value = 0.9961649491069209
$ws.Range("T4").Value = 0.9961649491069208
